$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -13.84609999999999
$ws.Range("E4").Value = 13.9276
$ws.Range("E5").Value = 13.0108
$ws.Range("C6").Value = -11.1705
$ws.Range("C7").Value = -11.669
$ws.Range("E8").Value = 14.1192
$ws.Range("C16").Value = -11.6002
$ws.Range("E16").Value = 12.9255
$ws.Range("C20").Value = -14.61260000000001
$ws.Range("E22").Value = 12.2543
